$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Dense matrix" comparison block (columns N:T) -----------------

# Header band (row 2) -- merged title over N2:S2
$ws.Range("N2:S2").Merge()
$ws.Range("N2").Value = "Dense matrix for the same design"
$ws.Range("N2:S2").Borders.LineStyle = 1
$ws.Range("N2:S2").HorizontalAlignment = -4108
$ws.Range("T2").Borders.LineStyle = 1

# Column headers (row 3)
$ws.Range("N3").Value = "numCols"
$ws.Range("O3").Value = "Resources"
$ws.Range("P3").Value = "percent"
$ws.Range("Q3").Value = "dataSize"
$ws.Range("R3").Value = "Time (ms)"
$ws.Range("S3").Value = "Gflops"
$ws.Range("T3").Value = "Ratio"
$ws.Range("N3:T3").Font.Bold = $true
$ws.Range("N3:T3").Borders.LineStyle = 1

# Data rows 4-15
$ws.Range("N4").Value = 8192
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 0.01
$ws.Range("Q4").Value = 663552
$ws.Range("R4").Value = 0.165078
$ws.Range("S4").Value = 0.01211500000000000077
$ws.Range("T4").Formula = "=S4/S4"

$ws.Range("P5").Value = 0.1
$ws.Range("Q5").Value = 6709248
$ws.Range("R5").Value = 0.152341
$ws.Range("S5").Value = 0.01312800000000000078
$ws.Range("T5").Formula = "=S5/S5"

$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 67108864
$ws.Range("R6").Value = 0.15112
$ws.Range("S6").Value = 0.01323399999999999924
$ws.Range("T6").Formula = "=S6/S6"

$ws.Range("O7").Value = 2
$ws.Range("P7").Value = 0.01
$ws.Range("Q7").Value = 663552
$ws.Range("R7").Value = "nil"
$ws.Range("S7").Value = "nil"

$ws.Range("P8").Value = 0.1
$ws.Range("Q8").Value = 6709248
$ws.Range("R8").Value = "nil"
$ws.Range("S8").Value = "nil"

$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 67108864
$ws.Range("R9").Value = "nil"
$ws.Range("S9").Value = "nil"

$ws.Range("O10").Value = 4
$ws.Range("P10").Value = 0.01
$ws.Range("Q10").Value = 663552
$ws.Range("R10").Value = 0.08806500000000000439
$ws.Range("S10").Value = 0.02271000000000000102
$ws.Range("T10").Formula = "=S10/S4"

$ws.Range("P11").Value = 0.1
$ws.Range("Q11").Value = 6709248
$ws.Range("R11").Value = 0.04345700000000000257
$ws.Range("S11").Value = 0.04602300000000000141
$ws.Range("T11").Formula = "=S11/S5"

$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 67108864
$ws.Range("R12").Value = 0.03940500000000000252
$ws.Range("S12").Value = 0.0507550000000000015
$ws.Range("T12").Formula = "=S12/S6"

$ws.Range("O13").Value = 8
$ws.Range("P13").Value = 0.01
$ws.Range("Q13").Value = 663552
$ws.Range("R13").Value = 0.121377
$ws.Range("S13").Value = 0.01647799999999999959
$ws.Range("T13").Formula = "=S13/S4"

$ws.Range("P14").Value = 0.1
$ws.Range("Q14").Value = 6709248
$ws.Range("R14").Value = 0.02953799999999999829
$ws.Range("S14").Value = 0.06771000000000000629
$ws.Range("T14").Formula = "=S14/S5"

$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 67108864
$ws.Range("R15").Value = 0.02085199999999999901
$ws.Range("S15").Value = 0.09591600000000000126
$ws.Range("T15").Formula = "=S15/S6"

# Merges for the repeated numCols / numPipes values (done before
# formatting so the border/alignment pass below lands uniformly across
# every cell of the merged ranges)
$ws.Range("N4:N15").Merge()
$ws.Range("O4:O6").Merge()
$ws.Range("O7:O9").Merge()
$ws.Range("O10:O12").Merge()
$ws.Range("O13:O15").Merge()

# Formatting for the data block
$ws.Range("N4:O15").Borders.LineStyle = 1
$ws.Range("N4:O15").HorizontalAlignment = -4108
$ws.Range("N4:O15").VerticalAlignment = -4108

$ws.Range("P4:T15").Borders.LineStyle = 1
$ws.Range("P4:T15").HorizontalAlignment = -4152

# Column K width tweak
$ws.Columns.Item(11).ColumnWidth = 9

# Stray bordered cell far below the table
$ws.Range("I26").Borders.LineStyle = 1

$ws.Range("U8").Select()
